$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit is a cyclic re-shuffle of the data held in rows 25, 27, 28, 29, 30 and 31
# (two separate cycles: 25<->30, and 27->31->28->29->27), so that the per-row
# "identity" columns (A, B, D, E, F, G, H, Q, R) -- and, for the 25/30 pair, the
# extra K/L/M/N/AC columns -- move to their new rows while the rest of each row
# (which is identical across all six rows) stays put.

# --- Row 25 (becomes what row 30 used to hold) ---
$ws.Range("A25").Value = 111936854
$ws.Range("B25").Value = 56414
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 100049
$ws.Range("F25").Value = "Spillkråka"
$ws.Range("G25").Value = "Dryocopus martius"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
$ws.Range("K25").Value = ""
$ws.Range("L25").Value = ""
$ws.Range("M25").Value = ""
$ws.Range("N25").Value = ""
$ws.Range("Q25").Value = 450998.3386916541
$ws.Range("R25").Value = 7087288.958247212
$ws.Range("AC25").Value = ""

# --- Row 27 (becomes what row 29 used to hold) ---
$ws.Range("A27").Value = 111936864
$ws.Range("B27").Value = 89423
$ws.Range("D27").Value = "NT"
$ws.Range("E27").Value = 5432
$ws.Range("F27").Value = "Granticka"
$ws.Range("G27").Value = "Porodaedalea chrysoloma"
$ws.Range("H27").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q27").Value = 451094.1385684713
$ws.Range("R27").Value = 7087212.607717684

# --- Row 28 (becomes what row 31 used to hold) ---
$ws.Range("A28").Value = 111936860
$ws.Range("B28").Value = 89423
$ws.Range("D28").Value = "NT"
$ws.Range("E28").Value = 5432
$ws.Range("F28").Value = "Granticka"
$ws.Range("G28").Value = "Porodaedalea chrysoloma"
$ws.Range("H28").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q28").Value = 450975.0281813644
$ws.Range("R28").Value = 7086982.857739178

# --- Row 29 (becomes what row 28 used to hold) ---
$ws.Range("A29").Value = 111936789
$ws.Range("B29").Value = 90087
$ws.Range("D29").Value = "LC"
$ws.Range("E29").Value = 3298
$ws.Range("F29").Value = "Trådticka"
$ws.Range("G29").Value = "Climacocystis borealis"
$ws.Range("H29").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q29").Value = 450955.1314140605
$ws.Range("R29").Value = 7087063.751596102

# --- Row 30 (becomes what row 25 used to hold) ---
$ws.Range("A30").Value = 111936793
$ws.Range("B30").Value = 56398
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 100109
$ws.Range("F30").Value = "Tretåig hackspett"
$ws.Range("G30").Value = "Picoides tridactylus"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("Q30").Value = 451088.7179698629
$ws.Range("R30").Value = 7087232.506422138
$ws.Range("AC30").Value = "ringhack äldre"

# --- Row 31 (becomes what row 27 used to hold) ---
$ws.Range("A31").Value = 111936892
$ws.Range("B31").Value = 77515
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 6425
$ws.Range("F31").Value = "Garnlav"
$ws.Range("G31").Value = "Alectoria sarmentosa"
$ws.Range("H31").Value = "(Ach.) Ach."
$ws.Range("Q31").Value = 451172.0902361136
$ws.Range("R31").Value = 7086726.569319103
